# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("AJ2").Value = 0.0138339531749992
$ws.Range("E2").Value = 0.1883562720413256
$ws.Range("F2").Value = 0.02510516389178956
$ws.Range("G2").Value = 0.2806140951908258
$ws.Range("J2").Value = 0.01799627644658349
$ws.Range("K2").Value = 0.09188788860966107
$ws.Range("M2").Value = 0.01772962084126473
$ws.Range("N2").Value = 0.08651575188835429
$ws.Range("O2").Value = 0.007515628400145689
$ws.Range("R2").Value = 0.0679215628021367
$ws.Range("S2").Value = 0.009354707817189264
$ws.Range("T2").Value = 0.009409078027060314
$ws.Range("U2").Value = 0.06249699576466432
$ws.Range("W2").Value = 0.0009130098233818617
$ws.Range("X2").Value = 0.08538596395270243
$ws.Range("Y2").Value = 0.01451133998321967
$ws.Range("Z2").Value = 0.02045269134469629
$ws.Range("AI3").Value = 0.008915787724415331
$ws.Range("AJ3").Value = 0.02160668316878013
$ws.Range("E3").Value = 0.2810718858277978
$ws.Range("F3").Value = 0.09570109906500639
$ws.Range("G3").Value = 0.3081852082731207
$ws.Range("J3").Value = 0.0043345091647187
$ws.Range("K3").Value = 0.1147149655833565
$ws.Range("L3").Value = 0.00793617876372193
$ws.Range("M3").Value = 0.02655650898660795
$ws.Range("N3").Value = 0.01184477026260017
$ws.Range("Q3").Value = 0.002125668525927968
$ws.Range("R3").Value = 0.007929201533682935
$ws.Range("S3").Value = 0.02547915724157416
$ws.Range("U3").Value = 0.01472853565955751
$ws.Range("X3").Value = 0.06350268928220477
$ws.Range("Z3").Value = 0.005367150936926981
$ws.Range("AH4").Value = 0.000849403767639064
$ws.Range("AI4").Value = 0.01497005020782152
$ws.Range("AJ4").Value = 0.01685582266043188
$ws.Range("E4").Value = 0.2847639121176406
$ws.Range("F4").Value = 0.1407597575663971
$ws.Range("G4").Value = 0.2709178435558056
$ws.Range("K4").Value = 0.1238003093150366
$ws.Range("L4").Value = 0.01895901933885585
$ws.Range("M4").Value = 0.03177633361973534
$ws.Range("N4").Value = 0.004106824489686503
$ws.Range("Q4").Value = 0.003039826768261577
$ws.Range("S4").Value = 0.01598442663664308
$ws.Range("U4").Value = 0.004849594737302712
$ws.Range("X4").Value = 0.06836687521874248
$ws.Range("AI5").Value = 0.003305170023756752
$ws.Range("AJ5").Value = 0.01872443188206918
$ws.Range("E5").Value = 0.2016516135182773
$ws.Range("F5").Value = 0.04034648798539808
$ws.Range("G5").Value = 0.3233782875579669
$ws.Range("H5").Value = 0.007408307741336323
$ws.Range("J5").Value = 0.03172615759207197
$ws.Range("K5").Value = 0.1075433285689083
$ws.Range("L5").Value = 0.01045456113914201
$ws.Range("M5").Value = 0.03323315708005286
$ws.Range("N5").Value = 0.02547289721365328
$ws.Range("R5").Value = 0.02467773573741563
$ws.Range("S5").Value = 0.01150139106056939
$ws.Range("U5").Value = 0.03517353274975801
$ws.Range("X5").Value = 0.09222974245457435
$ws.Range("Y5").Value = 0.01130436783930669
$ws.Range("Z5").Value = 0.02186882985574325
$ws.Range("AH6").Value = 0.002053884637459715
$ws.Range("AI6").Value = 0.01812342740895694
$ws.Range("D6").Value = 0.005630125271338065
$ws.Range("E6").Value = 0.2212141198256587
$ws.Range("F6").Value = 0.2265284776502698
$ws.Range("G6").Value = 0.1807452827802943
$ws.Range("J6").Value = 0.005106432762393669
$ws.Range("K6").Value = 0.1973088537913883
$ws.Range("L6").Value = 0.01555486107933129
$ws.Range("M6").Value = 0.003096913791262934
$ws.Range("Q6").Value = 0.009917407698003351
$ws.Range("T6").Value = 0.01354812936084224
$ws.Range("W6").Value = 0.01871049013293005
$ws.Range("X6").Value = 0.0824615938098706

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("AA2").Value = 0.9861660468250008
$ws.Range("AB2").Value = 0.9861660468250008
$ws.Range("AC2").Value = 0.9861660468250008
$ws.Range("AD2").Value = 0.9861660468250008
$ws.Range("AE2").Value = 0.9861660468250008
$ws.Range("AF2").Value = 0.9861660468250008
$ws.Range("AG2").Value = 0.9861660468250008
$ws.Range("AH2").Value = 0.9861660468250008
$ws.Range("AI2").Value = 0.9861660468250008
$ws.Range("E2").Value = 0.1883562720413256
$ws.Range("F2").Value = 0.2134614359331152
$ws.Range("G2").Value = 0.4940755311239409
$ws.Range("H2").Value = 0.4940755311239409
$ws.Range("I2").Value = 0.4940755311239409
$ws.Range("J2").Value = 0.5120718075705244
$ws.Range("K2").Value = 0.6039596961801854
$ws.Range("L2").Value = 0.6039596961801854
$ws.Range("M2").Value = 0.6216893170214501
$ws.Range("N2").Value = 0.7082050689098044
$ws.Range("O2").Value = 0.71572069730995
$ws.Range("P2").Value = 0.71572069730995
$ws.Range("Q2").Value = 0.71572069730995
$ws.Range("R2").Value = 0.7836422601120867
$ws.Range("S2").Value = 0.792996967929276
$ws.Range("T2").Value = 0.8024060459563362
$ws.Range("U2").Value = 0.8649030417210005
$ws.Range("V2").Value = 0.8649030417210005
$ws.Range("W2").Value = 0.8658160515443823
$ws.Range("X2").Value = 0.9512020154970848
$ws.Range("Y2").Value = 0.9657133554803045
$ws.Range("Z2").Value = 0.9861660468250008
$ws.Range("AA3").Value = 0.9694775291068046
$ws.Range("AB3").Value = 0.9694775291068046
$ws.Range("AC3").Value = 0.9694775291068046
$ws.Range("AD3").Value = 0.9694775291068046
$ws.Range("AE3").Value = 0.9694775291068046
$ws.Range("AF3").Value = 0.9694775291068046
$ws.Range("AG3").Value = 0.9694775291068046
$ws.Range("AH3").Value = 0.9694775291068046
$ws.Range("AI3").Value = 0.97839331683122
$ws.Range("E3").Value = 0.2810718858277978
$ws.Range("F3").Value = 0.3767729848928042
$ws.Range("G3").Value = 0.6849581931659249
$ws.Range("H3").Value = 0.6849581931659249
$ws.Range("I3").Value = 0.6849581931659249
$ws.Range("J3").Value = 0.6892927023306435
$ws.Range("K3").Value = 0.8040076679140001
$ws.Range("L3").Value = 0.811943846677722
$ws.Range("M3").Value = 0.83850035566433
$ws.Range("N3").Value = 0.8503451259269302
$ws.Range("O3").Value = 0.8503451259269302
$ws.Range("P3").Value = 0.8503451259269302
$ws.Range("Q3").Value = 0.8524707944528582
$ws.Range("R3").Value = 0.8603999959865412
$ws.Range("S3").Value = 0.8858791532281153
$ws.Range("T3").Value = 0.8858791532281153
$ws.Range("U3").Value = 0.9006076888876728
$ws.Range("V3").Value = 0.9006076888876728
$ws.Range("W3").Value = 0.9006076888876728
$ws.Range("X3").Value = 0.9641103781698777
$ws.Range("Y3").Value = 0.9641103781698777
$ws.Range("Z3").Value = 0.9694775291068046
$ws.Range("AA4").Value = 0.9673247233641075
$ws.Range("AB4").Value = 0.9673247233641075
$ws.Range("AC4").Value = 0.9673247233641075
$ws.Range("AD4").Value = 0.9673247233641075
$ws.Range("AE4").Value = 0.9673247233641075
$ws.Range("AF4").Value = 0.9673247233641075
$ws.Range("AG4").Value = 0.9673247233641075
$ws.Range("AH4").Value = 0.9681741271317466
$ws.Range("AI4").Value = 0.983144177339568
$ws.Range("E4").Value = 0.2847639121176406
$ws.Range("F4").Value = 0.4255236696840377
$ws.Range("G4").Value = 0.6964415132398433
$ws.Range("H4").Value = 0.6964415132398433
$ws.Range("I4").Value = 0.6964415132398433
$ws.Range("J4").Value = 0.6964415132398433
$ws.Range("K4").Value = 0.8202418225548799
$ws.Range("L4").Value = 0.8392008418937358
$ws.Range("M4").Value = 0.8709771755134711
$ws.Range("N4").Value = 0.8750840000031576
$ws.Range("O4").Value = 0.8750840000031576
$ws.Range("P4").Value = 0.8750840000031576
$ws.Range("Q4").Value = 0.8781238267714191
$ws.Range("R4").Value = 0.8781238267714191
$ws.Range("S4").Value = 0.8941082534080622
$ws.Range("T4").Value = 0.8941082534080622
$ws.Range("U4").Value = 0.8989578481453649
$ws.Range("V4").Value = 0.8989578481453649
$ws.Range("W4").Value = 0.8989578481453649
$ws.Range("X4").Value = 0.9673247233641075
$ws.Range("Y4").Value = 0.9673247233641075
$ws.Range("Z4").Value = 0.9673247233641075
$ws.Range("AA5").Value = 0.9779703980941741
$ws.Range("AB5").Value = 0.9779703980941741
$ws.Range("AC5").Value = 0.9779703980941741
$ws.Range("AD5").Value = 0.9779703980941741
$ws.Range("AE5").Value = 0.9779703980941741
$ws.Range("AF5").Value = 0.9779703980941741
$ws.Range("AG5").Value = 0.9779703980941741
$ws.Range("AH5").Value = 0.9779703980941741
$ws.Range("AI5").Value = 0.9812755681179308
$ws.Range("AJ5").Value = 1
$ws.Range("E5").Value = 0.2016516135182773
$ws.Range("F5").Value = 0.2419981015036753
$ws.Range("G5").Value = 0.5653763890616422
$ws.Range("H5").Value = 0.5727846968029785
$ws.Range("I5").Value = 0.5727846968029785
$ws.Range("J5").Value = 0.6045108543950505
$ws.Range("K5").Value = 0.7120541829639587
$ws.Range("L5").Value = 0.7225087441031007
$ws.Range("M5").Value = 0.7557419011831535
$ws.Range("N5").Value = 0.7812147983968069
$ws.Range("O5").Value = 0.7812147983968069
$ws.Range("P5").Value = 0.7812147983968069
$ws.Range("Q5").Value = 0.7812147983968069
$ws.Range("R5").Value = 0.8058925341342225
$ws.Range("S5").Value = 0.8173939251947918
$ws.Range("T5").Value = 0.8173939251947918
$ws.Range("U5").Value = 0.8525674579445498
$ws.Range("V5").Value = 0.8525674579445498
$ws.Range("W5").Value = 0.8525674579445498
$ws.Range("X5").Value = 0.9447972003991242
$ws.Range("Y5").Value = 0.9561015682384308
$ws.Range("Z5").Value = 0.9779703980941741
$ws.Range("AA6").Value = 0.9798226879535833
$ws.Range("AB6").Value = 0.9798226879535833
$ws.Range("AC6").Value = 0.9798226879535833
$ws.Range("AD6").Value = 0.9798226879535833
$ws.Range("AE6").Value = 0.9798226879535833
$ws.Range("AF6").Value = 0.9798226879535833
$ws.Range("AG6").Value = 0.9798226879535833
$ws.Range("AH6").Value = 0.981876572591043
$ws.Range("AI6").Value = 0.9999999999999999
$ws.Range("AJ6").Value = 0.9999999999999999
$ws.Range("D6").Value = 0.005630125271338065
$ws.Range("E6").Value = 0.2268442450969968
$ws.Range("F6").Value = 0.4533727227472666
$ws.Range("G6").Value = 0.6341180055275609
$ws.Range("H6").Value = 0.6341180055275609
$ws.Range("I6").Value = 0.6341180055275609
$ws.Range("J6").Value = 0.6392244382899546
$ws.Range("K6").Value = 0.8365332920813429
$ws.Range("L6").Value = 0.8520881531606741
$ws.Range("M6").Value = 0.8551850669519371
$ws.Range("N6").Value = 0.8551850669519371
$ws.Range("O6").Value = 0.8551850669519371
$ws.Range("P6").Value = 0.8551850669519371
$ws.Range("Q6").Value = 0.8651024746499404
$ws.Range("R6").Value = 0.8651024746499404
$ws.Range("S6").Value = 0.8651024746499404
$ws.Range("T6").Value = 0.8786506040107827
$ws.Range("U6").Value = 0.8786506040107827
$ws.Range("V6").Value = 0.8786506040107827
$ws.Range("W6").Value = 0.8973610941437127
$ws.Range("X6").Value = 0.9798226879535833
$ws.Range("Y6").Value = 0.9798226879535833
$ws.Range("Z6").Value = 0.9798226879535833

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D2").Value = 9
$ws.Range("F2").Value = 0.5120718075705244
$ws.Range("G2").Value = 7
$ws.Range("F3").Value = 0.6849581931659249
$ws.Range("F4").Value = 0.6964415132398433
$ws.Range("F5").Value = 0.5653763890616422
$ws.Range("D6").Value = 6
$ws.Range("F6").Value = 0.6341180055275609
$ws.Range("G6").Value = 4

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 13
$ws.Range("F2").Value = 0.7082050689098044
$ws.Range("G2").Value = 11
$ws.Range("D3").Value = 10
$ws.Range("F3").Value = 0.8040076679140001
$ws.Range("G3").Value = 8
$ws.Range("D4").Value = 10
$ws.Range("F4").Value = 0.8202418225548799
$ws.Range("G4").Value = 8
$ws.Range("D5").Value = 10
$ws.Range("F5").Value = 0.7120541829639587
$ws.Range("G5").Value = 8
$ws.Range("D6").Value = 10
$ws.Range("F6").Value = 0.8365332920813429
$ws.Range("G6").Value = 8

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 19
$ws.Range("F2").Value = 0.8024060459563362
$ws.Range("G2").Value = 17
$ws.Range("D3").Value = 10
$ws.Range("F3").Value = 0.8040076679140001
$ws.Range("G3").Value = 8
$ws.Range("D4").Value = 10
$ws.Range("F4").Value = 0.8202418225548799
$ws.Range("G4").Value = 8
$ws.Range("D5").Value = 17
$ws.Range("F5").Value = 0.8058925341342225
$ws.Range("G5").Value = 15
$ws.Range("F6").Value = 0.8365332920813429

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 23
$ws.Range("F2").Value = 0.9512020154970848
$ws.Range("G2").Value = 21
$ws.Range("D3").Value = 20
$ws.Range("F3").Value = 0.9006076888876728
$ws.Range("G3").Value = 18
$ws.Range("D4").Value = 23
$ws.Range("F4").Value = 0.9673247233641075
$ws.Range("G4").Value = 21
$ws.Range("F5").Value = 0.9447972003991242
$ws.Range("D6").Value = 23
$ws.Range("F6").Value = 0.9798226879535833
$ws.Range("G6").Value = 21

Write-Output "Edit complete"